# Updated cryptos list on Fri Jul 12 19:18:48 UTC 2024 with GitHub Actions
# Refreshes price / volume(1h) figures (and, for rows 46-47, swaps the
# FirstDigitalUSD / VeChain ranking) to match the latest scrape.
# A leading apostrophe is used on every assigned value to force Excel to
# store it as literal text (matching the original inlineStr cells) instead
# of auto-converting numeric-looking strings such as "1.00" -> 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.178.30"
$ws.Range("E2").Value = "'  +0.91%  "
$ws.Range("D3").Value = "'3.136.52"
$ws.Range("E3").Value = "'  +0.16%  "
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'535.16"
$ws.Range("D6").Value = "'139.44"
$ws.Range("E6").Value = "'  +0.97%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'3.135.15"
$ws.Range("E8").Value = "'  +0.19%  "
$ws.Range("D9").Value = "'0.468"
$ws.Range("E9").Value = "'  +4.53%  "
$ws.Range("E10").Value = "'  +1.09%  "
$ws.Range("E11").Value = "'  +0.15%  "
$ws.Range("E12").Value = "'  +4.57%  "
$ws.Range("D13").Value = "'3.675.52"
$ws.Range("E13").Value = "'  +0.02%  "
$ws.Range("D14").Value = "'0.136"
$ws.Range("E14").Value = "'  +1.13%  "
$ws.Range("D15").Value = "'25.67"
$ws.Range("E15").Value = "'  +0.73%  "
$ws.Range("E16").Value = "'  +0.17%  "
$ws.Range("D17").Value = "'58.285.86"
$ws.Range("E17").Value = "'  +0.77%  "
$ws.Range("D18").Value = "'3.140.29"
$ws.Range("E18").Value = "'  -0.22%  "
$ws.Range("D19").Value = "'6.06"
$ws.Range("E19").Value = "'  +1.20%  "
$ws.Range("D20").Value = "'12.77"
$ws.Range("E20").Value = "'  +0.60%  "
$ws.Range("D21").Value = "'8.16"
$ws.Range("E21").Value = "'  +2.62%  "
$ws.Range("D22").Value = "'361.39"
$ws.Range("E22").Value = "'  +3.01%  "
$ws.Range("E23").Value = "'  +0.05%  "
$ws.Range("D24").Value = "'69.10"
$ws.Range("E24").Value = "'  +0.92%  "
$ws.Range("D25").Value = "'0.507"
$ws.Range("E25").Value = "'  +0.01%  "
$ws.Range("E26").Value = "'  -1.45%  "
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "'  +2.71%  "
$ws.Range("E28").Value = "'  -4.15%  "
$ws.Range("D29").Value = "'7.35"
$ws.Range("E29").Value = "'  -2.28%  "
$ws.Range("D30").Value = "'1.88"
$ws.Range("E30").Value = "'  +0.45%  "
$ws.Range("D31").Value = "'6.12"
$ws.Range("E31").Value = "'  -0.12%  "
$ws.Range("D32").Value = "'21.46"
$ws.Range("E32").Value = "'  +1.53%  "
$ws.Range("E33").Value = "'  +2.20%  "
$ws.Range("D34").Value = "'1.15"
$ws.Range("E34").Value = "'  -2.49%  "
$ws.Range("D35").Value = "'158.61"
$ws.Range("E35").Value = "'  +0.25%  "
$ws.Range("D36").Value = "'6.09"
$ws.Range("E36").Value = "'  -1.56%  "
$ws.Range("D37").Value = "'25.88"
$ws.Range("E37").Value = "'  -1.65%  "
$ws.Range("E38").Value = "'  +1.61%  "
$ws.Range("D39").Value = "'1.67"
$ws.Range("E39").Value = "'  +2.99%  "
$ws.Range("D40").Value = "'0.0674"
$ws.Range("E40").Value = "'  +0.55%  "
$ws.Range("D41").Value = "'2.509.33"
$ws.Range("E41").Value = "'  +6.89%  "
$ws.Range("D42").Value = "'0.702"
$ws.Range("E42").Value = "'  -0.09%  "
$ws.Range("E43").Value = "'  -4.44%  "
$ws.Range("E44").Value = "'  +3.09%  "
$ws.Range("D45").Value = "'3.179.60"
$ws.Range("B46").Value = "'FirstDigitalUSD"
$ws.Range("C46").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "'  -0.03%  "
$ws.Range("B47").Value = "'VeChain"
$ws.Range("C47").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0269"
$ws.Range("E47").Value = "'  -0.93%  "
$ws.Range("D48").Value = "'0.996"
$ws.Range("D49").Value = "'6.09"
$ws.Range("E49").Value = "'  +0.96%  "
$ws.Range("D50").Value = "'19.93"
$ws.Range("D51").Value = "'0.743"
$ws.Range("E51").Value = "'  -3.35%  "
